$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 19:06"

# Row 4
$ws.Range("B4").Value = 6222086
$ws.Range("C4").Value = 10290
$ws.Range("D4").Value = 3460325
$ws.Range("E4").Value = 2573789
$ws.Range("G4").Value = 236
$ws.Range("H4").Value = 187972

# Row 5
$ws.Range("B5").Value = 3919452
$ws.Range("C5").Value = 8551
$ws.Range("E5").Value = 699991
$ws.Range("G5").Value = 212
$ws.Range("H5").Value = 121727

# Row 6
$ws.Range("B6").Value = 3758705
$ws.Range("C6").Value = 70766
$ws.Range("D6").Value = 2895175
$ws.Range("E6").Value = 797111
$ws.Range("G6").Value = 984
$ws.Range("H6").Value = 66419

# Row 21
$ws.Range("B21").Value = 271705
$ws.Range("C21").Value = 1572
$ws.Range("D21").Value = 245929
$ws.Range("E21").Value = 19359
$ws.Range("G21").Value = 47
$ws.Range("H21").Value = 6417

# Row 22
$ws.Range("B22").Value = 270189
$ws.Range("C22").Value = 978
$ws.Range("D22").Value = 207944
$ws.Range("E22").Value = 26754
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 35491

# Row 27
$ws.Range("B27").Value = 129182
$ws.Range("C27").Value = 234
$ws.Range("D27").Value = 114396
$ws.Range("E27").Value = 5657
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 9129

# Row 30
$ws.Range("B30").Value = 118122
$ws.Range("C30").Value = 1526
$ws.Range("D30").Value = 96206
$ws.Range("E30").Value = 20960
$ws.Range("G30").Value = 17
$ws.Range("H30").Value = 956

# Row 32
$ws.Range("B32").Value = 114309
$ws.Range("C32").Value = 542
$ws.Range("D32").Value = 101723
$ws.Range("E32").Value = 6015
$ws.Range("G32").Value = 15
$ws.Range("H32").Value = 6571

# Row 33
$ws.Range("E33").Value = 7385
$ws.Range("G33").Value = 65
$ws.Range("H33").Value = 1588

# Row 54
$ws.Range("B54").Value = 53304
$ws.Range("C54").Value = 1173
$ws.Range("D54").Value = 19487
$ws.Range("E54").Value = 32989
$ws.Range("G54").Value = 19
$ws.Range("H54").Value = 828

# Row 74
$ws.Range("B74").Value = 24832
$ws.Range("C74").Value = 214
$ws.Range("D74").Value = 18088
$ws.Range("E74").Value = 6319
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 425

# Row 78
$ws.Range("B78").Value = 19409
$ws.Range("C78").Value = 267
$ws.Range("E78").Value = 1344
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 414

# Row 80
$ws.Range("A80").Value = "Libano"
$ws.Range("B80").Value = 17777
$ws.Range("C80").Value = 469
$ws.Range("D80").Value = 4988
$ws.Range("E80").Value = 12618
$ws.Range("G80").Value = 4
$ws.Range("H80").Value = 171

# Row 81
$ws.Range("A81").Value = "Paraguay"
$ws.Range("B81").Value = 17662
$ws.Range("D81").Value = 9421
$ws.Range("E81").Value = 7915
$ws.Range("H81").Value = 326

# Row 117
$ws.Range("A117").Value = "Mozambique"
$ws.Range("B117").Value = 4039
$ws.Range("C117").Value = 123
$ws.Range("D117").Value = 2170
$ws.Range("E117").Value = 1846
$ws.Range("H117").Value = 23

# Row 118
$ws.Range("A118").Value = "Surinam"
$ws.Range("B118").Value = 4034
$ws.Range("D118").Value = 3140
$ws.Range("E118").Value = 823
$ws.Range("H118").Value = 71

# Row 119
$ws.Range("A119").Value = "Cuba"
$ws.Range("B119").Value = 4032
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 3378
$ws.Range("E119").Value = 560
$ws.Range("H119").Value = 94

# Row 120
$ws.Range("A120").Value = "Eslovaquia"
$ws.Range("B120").Value = 3989
$ws.Range("C120").Value = 72
$ws.Range("D120").Value = 2478
$ws.Range("E120").Value = 1478
$ws.Range("H120").Value = 33

# Row 121
$ws.Range("A121").Value = "Congo"
$ws.Range("B121").Value = 3979
$ws.Range("D121").Value = 1742
$ws.Range("E121").Value = 2159
$ws.Range("H121").Value = 78

# Row 127
$ws.Range("B127").Value = 3092
$ws.Range("C127").Value = 43
$ws.Range("E127").Value = 201
